$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.749.53"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "3.323.75"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.76"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.23"
$ws.Range("E6").Value = "  -3.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").Value = "3.320.42"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.37"
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "659.65"
$ws.Range("E14").Value = "  +4.69%  "
$ws.Range("D15").Value = "3.868.49"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").Value = "67.705.33"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").Value = "3.323.95"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.32"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.96"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("E23").Value = "  +5.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.95"
$ws.Range("E24").Value = "  -3.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.45"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("E26").Value = "  -3.51%  "
$ws.Range("E27").Value = "  -4.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.24"
$ws.Range("E28").Value = "  -3.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.37"
$ws.Range("E29").Value = "  +2.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.41"
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.15"
$ws.Range("E31").Value = "  +6.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "567.61"
$ws.Range("E32").Value = "  -5.28%  "
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.664.41"
$ws.Range("E36").Value = "  -7.02%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.45"
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.27"
$ws.Range("E38").Value = "  -6.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.21"
$ws.Range("E39").Value = "  +4.82%  "
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("E41").Value = "  -2.29%  "
$ws.Range("E42").Value = "  -4.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.35"
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.332"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("E45").Value = "  -3.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0406"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.60"
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.85"
$ws.Range("E51").Value = "  -0.82%  "
